# Add new "shen" dataset rows (261-273) to Sheet1, extend the used range /
# AutoFilter / _FilterDatabase defined name accordingly, add the trailing
# blank formatted rows (274-281), and move the active selection the way the
# authoring session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. New data rows 261-273
# ---------------------------------------------------------------------
$rows = @(
    @{ A="shen"; B="Proximal Tubule Cell"; C="LRP2, CUBN, SLC13A1"; E="PTC" },
    @{ A="shen"; B="Thin Limb Cell"; C="CRYAB, TACSTD2, SLC44A5, KLRG2, COL26A1, BOC"; E="TLC" },
    @{ A="shen"; B="Thick Ascending Limb Cell"; C="CASR, SLC12A1, UMOD"; E="TALC" },
    @{ A="shen"; B="Distal Convoluted Tubule Cell"; C="SLC12A3, CNNM2, FGF13, KLHL3, LHX1, TRPM6"; E="DCTC" },
    @{ A="shen"; B="Connecting Tubule Cell"; C="SLC8A1, SCN2A, HSD11B2, CALB1"; E="CTC" },
    @{ A="shen"; B="Principal Cell"; C="GATA3, AQP2, AQP3"; E="PC" },
    @{ A="shen"; B="Papillary Epithelial Cell"; C="TACSTD2, TP63, GPX2, FXYD3, KRT5"; E="PapEC" },
    @{ A="shen"; B="Intercalated Cell"; C="ATP6V0D2, ATP6V1C2, TMEM213, CLNK"; E="IC" },
    @{ A="shen"; B="Endothelial Cell"; C="CD34, PECAM1, PTPRB, MEIS2, FLT1, EMCN"; E="EC" },
    @{ A="shen"; B="Vascular Smooth Muscle Cell / Pericyte"; C="NOTCH3, PDGFRB, ITGA8"; E="VSMC/P" },
    @{ A="shen"; B="Fibroblast"; C="COL1A1, COL1A2, C7, NEGR1, FBLN5, DCN, CDH11"; E="FIB" },
    @{ A="shen"; B="Immune Cell"; C="PTPRC"; E="IMMC" },
    @{ A="shen"; B="Macrophage"; C="CD68, CD163"; E="M" }
)

$r = 261
foreach ($row in $rows) {
    $rng = $ws.Range("A" + $r + ":E" + $r)
    $rng.HorizontalAlignment = $xlCenter
    $rng.VerticalAlignment = $xlCenter

    $ws.Range("A" + $r).Value = $row.A
    $ws.Range("B" + $r).Value = $row.B
    $ws.Range("C" + $r).Value = $row.C
    $ws.Range("E" + $r).Value = $row.E
    # column D is left blank (formatted only), same as every other data row

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Trailing formatted-but-empty rows 274-281 (column A only)
# ---------------------------------------------------------------------
$blank = $ws.Range("A274:A281")
$blank.HorizontalAlignment = $xlCenter
$blank.VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 3. Grow the AutoFilter range from A1:E260 to A1:E273
#    (toggle off/on so the stored range actually updates)
# ---------------------------------------------------------------------
$full = $ws.Range("A1:E273")
$full.AutoFilter()
$full.AutoFilter()

# ---------------------------------------------------------------------
# 4. Update the hidden _FilterDatabase defined name to match
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$273"
    }
}

# ---------------------------------------------------------------------
# 5. Leave the selection where the authoring session left it
# ---------------------------------------------------------------------
$ws.Range("A261:A273").Select()
